$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (shared strings content changed per diff).
# Order matters here: new shared-string table entries are appended in the
# order cells are written, so match the author's edit order (E, F, D, A, B).
$ws.Range("E2").Value = "testeo "
$ws.Range("F2").Value = "tester"
$ws.Range("D2").Value = "FRIGORÍFICO"
$ws.Range("A2").Value = "505050abc"
$ws.Range("B2").Value = "testeo 5&*/&& visitante"

# F2 loses its underline formatting (style xf/font for underline is removed entirely)
$ws.Range("F2").Font.Underline = 0

# G2 (placa) is cleared out entirely in the new layout
$ws.Range("G2").ClearContents()

# Dates shift forward
$ws.Range("H2").Value = 44927
$ws.Range("I2").Value = 44957

# Update the active selection to match the saved view
$ws.Range("F7").Select()
